# update doc & resize format login page
#
# - Fix the "upload profile picture" route to use kebab-case.
# - Merge the standalone "teacher" role row into the upload-profile-picture
#   row's role cell (now "admin, teacher (only owner)"), matching the
#   style already used by D5/D6.
# - Fix the "classesByUser" route to use kebab-case.
# - Shrink row 8's height a bit (login-page table resize).
# - Move the active selection from F14 to D14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- content fixes -------------------------------------------------------

# /api/users/:userid/uploadProfilePicture -> .../upload-profile-picture
$ws.Range("B8").Value = "/api/users/:userid/upload-profile-picture"

# D8 used to just say "teacher"; pick up the same formatting D5/D6 use
# (the "admin, teacher (only owner)" role cells) before overwriting the text,
# so the cell's style matches the rest of that role column exactly.
$ws.Range("D5").Copy()
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("D8").Value = "admin, teacher (only owner)"

# /api/classesByUser -> /api/classes-by-user
$ws.Range("B19").Value = "/api/classes-by-user"

# --- formatting / layout --------------------------------------------------

# resize the login row (row 8 holds the upload-profile-picture / POST row)
$ws.Rows.Item(8).RowHeight = 13.8

# move the saved selection to D14
$ws.Range("D14").Select()
